# Update "想去人数" (interest count) figures in both the "展览" and
# "全部类型" sheets to the newly refreshed numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13
$ws1.Range("F3").Value = 1336
$ws1.Range("F8").Value = 11564
$ws1.Range("F9").Value = 4360
$ws1.Range("F12").Value = 22
$ws1.Range("F14").Value = 2539
$ws1.Range("F15").Value = 1084
$ws1.Range("F18").Value = 3634
$ws1.Range("F19").Value = 181
$ws1.Range("F21").Value = 11316
$ws1.Range("F22").Value = 11231

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13
$ws4.Range("F3").Value = 1336
$ws4.Range("F8").Value = 11564
$ws4.Range("F9").Value = 4360
$ws4.Range("F12").Value = 22
$ws4.Range("F14").Value = 2539
$ws4.Range("F16").Value = 1084
$ws4.Range("F19").Value = 3634
$ws4.Range("F20").Value = 181
$ws4.Range("F22").Value = 11316
$ws4.Range("F23").Value = 11231
